$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.750.68'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.874.43'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '325.21'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4585'
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3870'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07852'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9917'
$ws.Range('E10').Value = '  +3.35%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.78'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.886.89'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.994'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.710'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06958'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '88.47'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '28.765.55'
$ws.Range('E21').Value = '  +2.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.281'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.04'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.133'
$ws.Range('E24').Value = '  +2.64%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.129.41'
$ws.Range('E25').Value = '  +2.43%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '153.18'
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.24'
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.792'
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.971'
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '119.07'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09309'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9181'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.299'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.339'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.05758'
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.155'
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02070'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.697'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5638'
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1788'
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '9.943'
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.07214'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '11.81'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5298'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.155'
$ws.Range('E46').Value = '  +3.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.123'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '113.65'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.827'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.411'
$ws.Range('E50').Value = '  +3.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.004'
$ws.Range('E51').Value = '  +0.28%  '
